$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) and "全部类型" (index 4) both contain the same
# two data rows that need their "想去人数" (F column) counts bumped.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 255
    $ws.Range("F3").Value = 369
}
